# Applies the "output generated at 456a3b4" update to the 合肥-漫展信息 workbook.
#
# Changes needed on both the "展览" sheet (sheet 1) and the "全部类型" sheet
# (sheet 4):
#   1. Bump a handful of "想去人数" (interest-count) values in F2, F3, F4, F5, F17.
#   2. Insert a brand-new event row right before the existing "2024-06-22 Look
#      Look" row (old row 19), pushing every row at/after it down by one.
#   3. Fill the freshly inserted row with the new event's data.
#
# The "演出" and "本地生活" sheets are untouched by this revision.

$wb = $excel.ActiveWorkbook

function Update-InterestCounts($ws) {
    $ws.Range("F2").Value = 179
    $ws.Range("F3").Value = 428
    $ws.Range("F4").Value = 12472
    $ws.Range("F5").Value = 1277
    $ws.Range("F17").Value = 3979
}

function Insert-NewEventRow($ws, $lastOldRow) {
    # Row 19 becomes the new event; everything that used to live at 19+
    # shifts down to 20+ (values - including col A's running index - move
    # with their rows). The source data regenerates col A as a plain
    # sequential index (row - 1), so re-stamp it back to 18, 19, 20, ...
    # across every row from the insertion point through the old last row
    # (now one row further down), instead of leaving the shifted values.
    $ws.Rows("19:19").Insert()

    $ws.Range("A19").Value = 18
    # "2024-06-09" would otherwise auto-coerce to a real date; force text
    # for this one write only.
    $ws.Range("B19").NumberFormat = "@"
    $ws.Range("B19").Value = "2024-06-09"
    $ws.Range("C19").Value = "合肥·第六届环形宇宙动漫游戏嘉年华内场票·赵成晨"
    $ws.Range("D19").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
    $ws.Range("E19").Value = "2024.06.09 09:30-06.09 17:00"
    $ws.Range("F19").Value = 1
    $ws.Range("G19").Value = 238
    $ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=84863"
    $ws.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202404/I5S4Ih2M1714031127805.jpeg"

    # Match the formatting of the surrounding data rows: col A keeps the
    # bordered/bold index style, the rest stay plain (this also clears the
    # one-off "@" text format stamped onto B19 above, matching how every
    # other date-looking cell in the sheet is stored: plain style, text
    # value).
    $ws.Range("A18").Copy()
    $ws.Range("A19").PasteSpecial(-4122)
    $ws.Range("B18:I18").Copy()
    $ws.Range("B19:I19").PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    # Re-sequence the running index in col A for every row pushed down by
    # the insert. Col A's index is always (row number - 1) throughout this
    # table (row 2 => 1, row 19 => 18, ...), so simply restamp that pattern
    # over every row from just after the new row through the (shifted) old
    # last row.
    for ($r = 20; $r -le ($lastOldRow + 1); $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lastOldRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    Update-InterestCounts $ws
    Insert-NewEventRow $ws $lastOldRow
}
